$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

$ws.Range("D2").Value = 113722
$ws.Range("E2").Value = 1720
$ws.Range("F2").Value = 1720
$ws.Range("G2").Value = 126
$ws.Range("H2").Value = -138
$ws.Range("I2").Value = -132
$ws.Range("J2").Value = -6
$ws.Range("K2").Value = 48968
$ws.Range("L2").Value = 33316
$ws.Range("M2").Value = 15652
$ws.Range("N2").Value = 14267
$ws.Range("O2").Value = 1385
$ws.Range("P2").Value = 1938
$ws.Range("Q2").Value = 3484
$ws.Range("R2").Value = -230
$ws.Range("S2").Value = -2370
$ws.Range("T2").Value = 88
$ws.Range("U2").Value = 3395
$ws.Range("V2").Value = 15571
$ws.Range("W2").Value = 1.51
$ws.Range("X2").Value = -0.12
$ws.Range("Y2").Value = -0.91
$ws.Range("Z2").Value = -0.28
$ws.Range("AA2").Value = 212.85
$ws.Range("AB2").Value = 620.8
$ws.Range("AC2").Value = -340
$ws.Range("AD2").Value = -84.31999999999999
$ws.Range("AE2").Value = 36905
$ws.Range("AF2").Value = 0.78
$ws.Range("AG2").Value = 300
$ws.Range("AH2").Value = 1.05
$ws.Range("AI2").Value = -87.91
$ws.Range("AJ2").Value = 38760000
$ws.Range("D3").Value = 132245
$ws.Range("E3").Value = 817
$ws.Range("F3").Value = 817
$ws.Range("G3").Value = -2152
$ws.Range("H3").Value = -2171
$ws.Range("I3").Value = -2354
$ws.Range("J3").Value = 183
$ws.Range("K3").Value = 53835
$ws.Range("L3").Value = 38402
$ws.Range("M3").Value = 15433
$ws.Range("N3").Value = 11993
$ws.Range("O3").Value = 3440
$ws.Range("P3").Value = 1938
$ws.Range("Q3").Value = 2977
$ws.Range("R3").Value = -2776
$ws.Range("S3").Value = 276
$ws.Range("T3").Value = 331
$ws.Range("U3").Value = 2646
$ws.Range("V3").Value = 17577
$ws.Range("W3").Value = 0.62
$ws.Range("X3").Value = -1.64
$ws.Range("Y3").Value = -17.93
$ws.Range("Z3").Value = -4.22
$ws.Range("AA3").Value = 248.82
$ws.Range("AB3").Value = 494.98
$ws.Range("AC3").Value = -6074
$ws.Range("AD3").Value = -5.56
$ws.Range("AE3").Value = 31023
$ws.Range("AF3").Value = 1.09
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 0.59
$ws.Range("AI3").Value = -3.28
$ws.Range("AJ3").Value = 38760000
$ws.Range("D4").Value = 119667
$ws.Range("E4").Value = 1741
$ws.Range("F4").Value = 1741
$ws.Range("G4").Value = 1332
$ws.Range("H4").Value = 848
$ws.Range("I4").Value = 581
$ws.Range("J4").Value = 266
$ws.Range("K4").Value = 51773
$ws.Range("L4").Value = 35319
$ws.Range("M4").Value = 16454
$ws.Range("N4").Value = 12489
$ws.Range("O4").Value = 3965
$ws.Range("P4").Value = 1938
$ws.Range("Q4").Value = 2303
$ws.Range("R4").Value = -2180
$ws.Range("S4").Value = -2377
$ws.Range("T4").Value = 779
$ws.Range("U4").Value = 1524
$ws.Range("V4").Value = 15280
$ws.Range("W4").Value = 1.46
$ws.Range("X4").Value = 0.71
$ws.Range("Y4").Value = 4.75
$ws.Range("Z4").Value = 1.6
$ws.Range("AA4").Value = 214.66
$ws.Range("AB4").Value = 521.8200000000001
$ws.Range("AC4").Value = 1500
$ws.Range("AD4").Value = 19.5
$ws.Range("AE4").Value = 32305
$ws.Range("AF4").Value = 0.91
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 0.85
$ws.Range("AI4").Value = 16.62
$ws.Range("AJ4").Value = 38760000
$ws.Range("D5").Value = 128272
$ws.Range("E5").Value = 2123
$ws.Range("F5").Value = 2123
$ws.Range("G5").Value = 1779
$ws.Range("H5").Value = 882
$ws.Range("I5").Value = 602
$ws.Range("J5").Value = 280
$ws.Range("K5").Value = 49676
$ws.Range("L5").Value = 33821
$ws.Range("M5").Value = 15855
$ws.Range("N5").Value = 11944
$ws.Range("O5").Value = 3911
$ws.Range("P5").Value = 1938
$ws.Range("Q5").Value = 1806
$ws.Range("R5").Value = -102
$ws.Range("S5").Value = -669
$ws.Range("T5").Value = 744
$ws.Range("U5").Value = 1063
$ws.Range("V5").Value = 14066
$ws.Range("W5").Value = 1.65
$ws.Range("X5").Value = 0.6899999999999999
$ws.Range("Y5").Value = 4.93
$ws.Range("Z5").Value = 1.74
$ws.Range("AA5").Value = 213.32
$ws.Range("AB5").Value = 547.95
$ws.Range("AC5").Value = 1553
$ws.Range("AD5").Value = 17.48
$ws.Range("AE5").Value = 30895
$ws.Range("AF5").Value = 0.88
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 0.92
$ws.Range("AI5").Value = 16.05
$ws.Range("AJ5").Value = 38760000
$ws.Range("D6").Value = 99882
$ws.Range("E6").Value = 1657
$ws.Range("F6").Value = 1657
$ws.Range("G6").Value = 1566
$ws.Range("H6").Value = -360
$ws.Range("I6").Value = -704
$ws.Range("K6").Value = 50381
$ws.Range("L6").Value = 34885
$ws.Range("M6").Value = 15495
$ws.Range("N6").Value = 11285
$ws.Range("P6").Value = 1938
$ws.Range("Q6").Value = 308
$ws.Range("R6").Value = -693
$ws.Range("S6").Value = -73
$ws.Range("T6").Value = 1030
$ws.Range("U6").Value = -723
$ws.Range("V6").Value = 14307
$ws.Range("W6").Value = 1.66
$ws.Range("X6").Value = -0.36
$ws.Range("Y6").Value = -6.06
$ws.Range("Z6").Value = -0.72
$ws.Range("AA6").Value = 225.13
$ws.Range("AB6").Value = 506.35
$ws.Range("AC6").Value = -1815
$ws.Range("AD6").Value = -8.460000000000001
$ws.Range("AE6").Value = 29191
$ws.Range("AF6").Value = 0.53
$ws.Range("AG6").Value = 250
$ws.Range("AH6").Value = 1.63
$ws.Range("AI6").Value = -13.74
$ws.Range("AJ6").Value = 38760000
$ws.Range("D7").Value = 105783
$ws.Range("E7").Value = 1666
$ws.Range("G7").Value = 2230
$ws.Range("H7").Value = 2179
$ws.Range("I7").Value = 1870
$ws.Range("K7").Value = 54201
$ws.Range("L7").Value = 36270
$ws.Range("M7").Value = 17930
$ws.Range("N7").Value = 13406
$ws.Range("P7").Value = 1940
$ws.Range("Q7").Value = 2281
$ws.Range("R7").Value = -2
$ws.Range("S7").Value = -917
$ws.Range("T7").Value = 913
$ws.Range("U7").Value = 1665
$ws.Range("W7").Value = 1.57
$ws.Range("X7").Value = 2.06
$ws.Range("Y7").Value = 15.15
$ws.Range("Z7").Value = 4.17
$ws.Range("AA7").Value = 202.28
$ws.Range("AC7").Value = 4826
$ws.Range("AD7").Value = 2.86
$ws.Range("AE7").Value = 34677
$ws.Range("AF7").Value = 0.4
$ws.Range("AG7").Value = 250
$ws.Range("AH7").Value = 1.81
$ws.Range("AI7").Value = 5.18
$ws.Range("D8").Value = 108328
$ws.Range("E8").Value = 1719
$ws.Range("G8").Value = 1509
$ws.Range("H8").Value = 1025
$ws.Range("I8").Value = 762
$ws.Range("K8").Value = 54165
$ws.Range("L8").Value = 36964
$ws.Range("M8").Value = 17200
$ws.Range("N8").Value = 12853
$ws.Range("P8").Value = 1940
$ws.Range("Q8").Value = 2158
$ws.Range("R8").Value = -1123
$ws.Range("S8").Value = -364
$ws.Range("T8").Value = 913
$ws.Range("U8").Value = 1235
$ws.Range("W8").Value = 1.59
$ws.Range("X8").Value = 0.95
$ws.Range("Y8").Value = 5.83
$ws.Range("Z8").Value = 1.91
$ws.Range("AA8").Value = 214.9
$ws.Range("AC8").Value = 1967
$ws.Range("AD8").Value = 6.15
$ws.Range("AE8").Value = 33246
$ws.Range("AF8").Value = 0.36
$ws.Range("AG8").Value = 275
$ws.Range("AH8").Value = 2.27
$ws.Range("AI8").Value = 13.98
$ws.Range("D9").Value = 112332
$ws.Range("E9").Value = 1914
$ws.Range("G9").Value = 1726
$ws.Range("H9").Value = 1168
$ws.Range("I9").Value = 874
$ws.Range("K9").Value = 55778
$ws.Range("L9").Value = 37496
$ws.Range("M9").Value = 18282
$ws.Range("N9").Value = 13702
$ws.Range("P9").Value = 1940
$ws.Range("Q9").Value = 1927
$ws.Range("R9").Value = -1288
$ws.Range("S9").Value = -279
$ws.Range("T9").Value = 1063
$ws.Range("U9").Value = 1166
$ws.Range("W9").Value = 1.7
$ws.Range("X9").Value = 1.04
$ws.Range("Y9").Value = 6.58
$ws.Range("Z9").Value = 2.12
$ws.Range("AA9").Value = 205.1
$ws.Range("AC9").Value = 2255
$ws.Range("AD9").Value = 5.37
$ws.Range("AE9").Value = 35444
$ws.Range("AF9").Value = 0.34
$ws.Range("AG9").Value = 275
$ws.Range("AH9").Value = 2.27
$ws.Range("AI9").Value = 12.2
